# Updated K (column G) values for rows 2-60, per regen of save_data
# using K instead of Strike# (std/mean recalculated and s_vals rewritten).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 0
    3 = 2
    4 = 2
    5 = 0
    6 = 1
    7 = 2
    9 = 0
    10 = 0
    12 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 2
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 2
    30 = 1
    31 = 0
    32 = 1
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 1
    38 = 0
    39 = 1
    40 = 2
    41 = 1
    42 = 0
    43 = 3
    44 = 0
    45 = 1
    46 = 1
    47 = 1
    48 = 2
    49 = 1
    50 = 1
    51 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 1
    57 = 3
    58 = 2
    59 = 1
    60 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
